$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(3)

# The text range already reads as "An image" (PowerPoint concatenates the
# underlying runs transparently), so assigning the same string back would be
# a no-op and wouldn't merge the three runs ("An", " ", "image") into one.
# Force a change first so the engine rewrites the paragraph into a single run.
$shape.TextFrame.TextRange.Text = "__tmp__"
$shape.TextFrame.TextRange.Text = "An image"
